# Updates cryptos list prices / volumes (and a few coin row re-ordering /
# replacements) to match the "Updated cryptos list" GitHub Actions commit.
# Numeric-looking price strings are prefixed with a literal leading
# single-quote so Excel stores them as text (matching the original
# inlineStr/text cells) instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.294.02'
$ws.Cells.Item(2, 5).Value = '  +0.71%  '
$ws.Cells.Item(3, 4).Value = '3.741.20'
$ws.Cells.Item(3, 5).Value = '  +0.09%  '
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).Value = '''593.09'
$ws.Cells.Item(5, 5).Value = '  -0.14%  '
$ws.Cells.Item(6, 4).Value = '''166.11'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '
$ws.Cells.Item(7, 4).Value = '3.739.84'
$ws.Cells.Item(7, 5).Value = '  +0.13%  '
$ws.Cells.Item(9, 4).Value = '''0.518'
$ws.Cells.Item(9, 5).Value = '  -0.09%  '
$ws.Cells.Item(10, 5).Value = '  -0.50%  '
$ws.Cells.Item(11, 4).Value = '''6.43'
$ws.Cells.Item(11, 5).Value = '  +0.03%  '
$ws.Cells.Item(12, 5).Value = '  +0.60%  '
$ws.Cells.Item(13, 4).Value = '''0.0000259'
$ws.Cells.Item(13, 5).Value = '  -2.83%  '
$ws.Cells.Item(14, 4).Value = '''36.17'
$ws.Cells.Item(14, 5).Value = '  +0.54%  '
$ws.Cells.Item(15, 4).Value = '4.368.46'
$ws.Cells.Item(15, 5).Value = '  +0.02%  '
$ws.Cells.Item(16, 4).Value = '3.742.79'
$ws.Cells.Item(16, 5).Value = '  +0.18%  '
$ws.Cells.Item(17, 4).Value = '68.291.71'
$ws.Cells.Item(17, 5).Value = '  +0.85%  '
$ws.Cells.Item(18, 4).Value = '''17.76'
$ws.Cells.Item(18, 5).Value = '  -3.12%  '
$ws.Cells.Item(19, 2).Value = 'Polkadot'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(19, 4).Value = '''6.98'
$ws.Cells.Item(19, 5).Value = '  -0.83%  '
$ws.Cells.Item(20, 2).Value = 'TRON'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(20, 4).Value = '''0.112'
$ws.Cells.Item(20, 5).Value = '  +0.53%  '
$ws.Cells.Item(21, 4).Value = '''10.67'
$ws.Cells.Item(21, 5).Value = '  +1.99%  '
$ws.Cells.Item(22, 4).Value = '''465.32'
$ws.Cells.Item(22, 5).Value = '  +0.38%  '
$ws.Cells.Item(23, 4).Value = '''0.694'
$ws.Cells.Item(23, 5).Value = '  -1.01%  '
$ws.Cells.Item(24, 4).Value = '''83.84'
$ws.Cells.Item(24, 5).Value = '  +1.16%  '
$ws.Cells.Item(25, 5).Value = '  +6.58%  '
$ws.Cells.Item(26, 5).Value = '  -0.92%  '
$ws.Cells.Item(27, 4).Value = '''11.87'
$ws.Cells.Item(27, 5).Value = '  -1.01%  '
$ws.Cells.Item(28, 4).Value = '''10.04'
$ws.Cells.Item(28, 5).Value = '  -2.49%  '
$ws.Cells.Item(29, 4).Value = '''0.999'
$ws.Cells.Item(29, 5).Value = '  -0.10%  '
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).Value = '''2.76'
$ws.Cells.Item(30, 5).Value = '  -4.19%  '
$ws.Cells.Item(31, 2).Value = 'NEARProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(31, 4).Value = '''7.29'
$ws.Cells.Item(31, 5).Value = '  -1.38%  '
$ws.Cells.Item(32, 2).Value = 'EthereumClassic'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(32, 4).Value = '''29.75'
$ws.Cells.Item(32, 5).Value = '  -0.39%  '
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).Value = '''2.16'
$ws.Cells.Item(33, 5).Value = '  -1.39%  '
$ws.Cells.Item(34, 2).Value = 'Aptos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(34, 4).Value = '''9.15'
$ws.Cells.Item(34, 5).Value = '  +1.07%  '
$ws.Cells.Item(35, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(35, 4).Value = '''0.999'
$ws.Cells.Item(35, 5).Value = '  --%  '
$ws.Cells.Item(36, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(36, 4).Value = '3.695.75'
$ws.Cells.Item(36, 5).Value = '  +0.20%  '
$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(37, 4).Value = '''0.100'
$ws.Cells.Item(37, 5).Value = '  -1.40%  '
$ws.Cells.Item(38, 2).Value = 'dogwifhat'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(38, 4).Value = '''3.41'
$ws.Cells.Item(38, 5).Value = '  -3.86%  '
$ws.Cells.Item(39, 2).Value = 'Kaspa'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(39, 4).Value = '''0.138'
$ws.Cells.Item(39, 5).Value = '  +0.59%  '
$ws.Cells.Item(40, 2).Value = 'Mantle'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(40, 4).Value = '''0.998'
$ws.Cells.Item(40, 5).Value = '  +0.51%  '
$ws.Cells.Item(41, 2).Value = 'Filecoin'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(41, 4).Value = '''5.77'
$ws.Cells.Item(41, 5).Value = '  +0.50%  '
$ws.Cells.Item(42, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(42, 4).Value = '''0.999'
$ws.Cells.Item(42, 5).Value = '  -0.11%  '
$ws.Cells.Item(43, 2).Value = 'USDe'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(43, 5).Value = '  -0.02%  '
$ws.Cells.Item(44, 2).Value = 'TheGraph'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(44, 4).Value = '''0.301'
$ws.Cells.Item(44, 5).Value = '  -1.99%  '
$ws.Cells.Item(45, 2).Value = 'Arweave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(45, 4).Value = '''43.37'
$ws.Cells.Item(45, 5).Value = '  +11.52%  '
$ws.Cells.Item(46, 2).Value = 'OKB'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(46, 4).Value = '''46.65'
$ws.Cells.Item(46, 5).Value = '  +3.43%  '
$ws.Cells.Item(47, 2).Value = 'Stacks'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(47, 4).Value = '''1.91'
$ws.Cells.Item(47, 5).Value = '  +0.06%  '
$ws.Cells.Item(48, 2).Value = 'Cosmos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(48, 4).Value = '''8.45'
$ws.Cells.Item(48, 5).Value = '  -0.98%  '
$ws.Cells.Item(49, 2).Value = 'Bittensor'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(49, 4).Value = '''390.36'
$ws.Cells.Item(49, 5).Value = '  -1.63%  '
$ws.Cells.Item(50, 2).Value = 'Monero'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(50, 4).Value = '''144.44'
$ws.Cells.Item(50, 5).Value = '  +0.37%  '
$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).Value = '2.745.40'
$ws.Cells.Item(51, 5).Value = '  +2.69%  '
